$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.245.55'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.603.43'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.22%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.39'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3770'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.01'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3645'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.278'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.00%  '

$ws.Range("B11").Value = 'BinanceUSD'
$ws.Range("C11").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08137'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.79'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.593'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.433'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001248'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.605.38'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.01'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06935'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.19'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.537'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.95'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.239.34'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.081'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.379'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.21'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.67'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.260'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.70'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.393'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.736'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.782.20'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9637'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07502'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.13%  '

$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.37'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.45%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02754'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2538'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.40%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08824'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.44%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.103'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.384'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7113'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.47'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.60'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6560'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.317'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9995'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.012'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.59'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07944'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.207'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.35%  '
